# Auto-applied data refresh for the Leve profit-tracking sheets.
# Mirrors a scheduled runner pulling fresh Market Board prices into the
# currentAveragePrice* / Leve profit columns (H:N) for specific rows.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 18
$ws.Range("H18").Value = 17835.75
$ws.Range("I18").Value = 23481
$ws.Range("K18").Value = 23481
$ws.Range("M18").Value = -23197
# Row 33
$ws.Range("H33").Value = 1947.7778
$ws.Range("I33").Value = 2132.8572
$ws.Range("K33").Value = 2132.8572
$ws.Range("M33").Value = -1903.8572
# Row 88
$ws.Range("H88").Value = 22263984
$ws.Range("I88").Value = 83336984
$ws.Range("J88").Value = 55621.816
$ws.Range("K88").Value = 83336984
$ws.Range("L88").Value = 55621.816
$ws.Range("M88").Value = -83336578
$ws.Range("N88").Value = -56433.816
# Row 91
$ws.Range("H91").Value = 22263984
$ws.Range("I91").Value = 83336984
$ws.Range("J91").Value = 55621.816
$ws.Range("K91").Value = 83336984
$ws.Range("L91").Value = 55621.816
$ws.Range("M91").Value = -83335580
$ws.Range("N91").Value = -58429.816
# Row 100
$ws.Range("H100").Value = 2438.7
$ws.Range("I100").Value = 1527.8182
$ws.Range("K100").Value = 1527.8182
$ws.Range("M100").Value = -986.8181999999999
# Row 135
$ws.Range("H135").Value = 3334633.2
$ws.Range("I135").Value = 3334633.2
$ws.Range("K135").Value = 30011698.8
$ws.Range("M135").Value = -30009163.8

$ws = $wb.Worksheets.Item("ARM")
# Row 76
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
# Row 79
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()
# Row 122
$ws.Range("H122").Value = 33999.715
$ws.Range("J122").Value = 7199.6
$ws.Range("L122").Value = 21598.8
$ws.Range("N122").Value = -26498.8

$ws = $wb.Worksheets.Item("BSM")
# Row 19
$ws.Range("H19").Value = 1300
$ws.Range("I19").Value = 1000
$ws.Range("J19").Value = 1500
$ws.Range("K19").Value = 1000
$ws.Range("L19").Value = 1500
$ws.Range("M19").Value = -827
$ws.Range("N19").Value = -1846
# Row 20
$ws.Range("H20").Value = 8773705
$ws.Range("I20").Value = 20835272
$ws.Range("K20").Value = 20835272
$ws.Range("M20").Value = -20835025
# Row 86
$ws.Range("H86").Value = 35752390
$ws.Range("I86").Value = 51813.85
$ws.Range("K86").Value = 51813.85
$ws.Range("M86").Value = -50690.85
# Row 89
$ws.Range("H89").Value = 35752390
$ws.Range("I89").Value = 51813.85
$ws.Range("K89").Value = 259069.25
$ws.Range("M89").Value = -253453.25
# Row 107
$ws.Range("I107").Value = 53574030
$ws.Range("J107").Value = 6179.6665
$ws.Range("K107").Value = 53574030
$ws.Range("L107").Value = 6179.6665
$ws.Range("M107").Value = -53572110
$ws.Range("N107").Value = -10019.6665

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 4139
$ws.Range("I16").Value = 2084.3333
$ws.Range("J16").Value = 6380.4546
$ws.Range("K16").Value = 2084.3333
$ws.Range("L16").Value = 6380.4546
$ws.Range("M16").Value = -1797.3333
$ws.Range("N16").Value = -6954.4546
# Row 31
$ws.Range("H31").Value = 6610.654
$ws.Range("I31").Value = 2707.5356
$ws.Range("J31").Value = 11164.292
$ws.Range("K31").Value = 2707.5356
$ws.Range("L31").Value = 11164.292
$ws.Range("M31").Value = -2412.5356
$ws.Range("N31").Value = -11754.292
# Row 34
$ws.Range("H34").Value = 6610.654
$ws.Range("I34").Value = 2707.5356
$ws.Range("J34").Value = 11164.292
$ws.Range("K34").Value = 2707.5356
$ws.Range("L34").Value = 11164.292
$ws.Range("M34").Value = -2505.5356
$ws.Range("N34").Value = -11568.292
# Row 58
$ws.Range("H58").Value = 11116581
$ws.Range("I58").Value = 17242956
$ws.Range("K58").Value = 17242956
$ws.Range("M58").Value = -17242753
# Row 62
$ws.Range("H62").Value = 8866.75
$ws.Range("J62").Value = 7005.3335
$ws.Range("L62").Value = 7005.3335
$ws.Range("N62").Value = -8253.333500000001
# Row 65
$ws.Range("H65").Value = 8866.75
$ws.Range("J65").Value = 7005.3335
$ws.Range("L65").Value = 35026.6675
$ws.Range("N65").Value = -41266.6675
# Row 86
$ws.Range("H86").Value = 9005420
$ws.Range("J86").Value = 129625
$ws.Range("L86").Value = 129625
$ws.Range("N86").Value = -131871
# Row 89
$ws.Range("H89").Value = 9005420
$ws.Range("J89").Value = 129625
$ws.Range("L89").Value = 648125
$ws.Range("N89").Value = -659357
# Row 113
$ws.Range("H113").Value = 4139
$ws.Range("I113").Value = 2084.3333
$ws.Range("J113").Value = 6380.4546
$ws.Range("K113").Value = 2084.3333
$ws.Range("L113").Value = 6380.4546
$ws.Range("M113").Value = 85.66670000000022
$ws.Range("N113").Value = -10720.4546
# Row 134
$ws.Range("H134").Value = 8040.1035
$ws.Range("I134").Value = 4066.8333
$ws.Range("J134").Value = 9076.608
$ws.Range("K134").Value = 12200.4999
$ws.Range("L134").Value = 27229.824
$ws.Range("M134").Value = -9665.499899999999
$ws.Range("N134").Value = -32299.824
# Row 136
$ws.Range("H136").Value = 11116581
$ws.Range("I136").Value = 17242956
$ws.Range("K136").Value = 51728868
$ws.Range("M136").Value = -51726318

$ws = $wb.Worksheets.Item("CUL")
# Row 14
$ws.Range("H14").Value = 13889236
$ws.Range("I14").Value = 13889236
$ws.Range("K14").Value = 41667708
$ws.Range("M14").Value = -41667535
# Row 44
$ws.Range("H44").Value = 650.8182
$ws.Range("I44").Value = 212
$ws.Range("J44").Value = 1016.5
$ws.Range("K44").Value = 636
$ws.Range("L44").Value = 3049.5
$ws.Range("M44").Value = -238
$ws.Range("N44").Value = -3845.5
# Row 61
$ws.Range("H61").Value = 332.72726
$ws.Range("I61").Value = 88.75
$ws.Range("J61").Value = 983.3333
$ws.Range("K61").Value = 266.25
$ws.Range("L61").Value = 2949.9999
$ws.Range("M61").Value = -51.25
$ws.Range("N61").Value = -3379.9999
# Row 88
$ws.Range("H88").Value = 2500
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()
# Row 91
$ws.Range("H91").Value = 2500
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()
# Row 98
$ws.Range("H98").Value = 3775.818
$ws.Range("I98").Value = 2995.5
$ws.Range("J98").Value = 4221.7144
$ws.Range("K98").Value = 8986.5
$ws.Range("L98").Value = 12665.1432
$ws.Range("M98").Value = -7488.5
$ws.Range("N98").Value = -15661.1432

$ws = $wb.Worksheets.Item("GSM")
# Row 35
$ws.Range("H35").Value = 21780.5
$ws.Range("I35").Value = 13561
$ws.Range("J35").Value = 30000
$ws.Range("K35").Value = 13561
$ws.Range("L35").Value = 30000
$ws.Range("M35").Value = -13263
$ws.Range("N35").Value = -30596
# Row 80
$ws.Range("H80").Value = 5191.3335
$ws.Range("I80").Value = 6024.5
$ws.Range("J80").Value = 4774.75
$ws.Range("K80").Value = 6024.5
$ws.Range("L80").Value = 4774.75
$ws.Range("M80").Value = -5026.5
$ws.Range("N80").Value = -6770.75
# Row 83
$ws.Range("H83").Value = 5191.3335
$ws.Range("I83").Value = 6024.5
$ws.Range("J83").Value = 4774.75
$ws.Range("K83").Value = 30122.5
$ws.Range("L83").Value = 23873.75
$ws.Range("M83").Value = -25130.5
$ws.Range("N83").Value = -33857.75
# Row 113
$ws.Range("H113").Value = 7108.9688
$ws.Range("I113").Value = 3254.4546
$ws.Range("K113").Value = 3254.4546
$ws.Range("M113").Value = -1084.4546
# Row 126
$ws.Range("H126").Value = 3216.5334
$ws.Range("I126").Value = 2579.1667
$ws.Range("K126").Value = 7737.500100000001
$ws.Range("M126").Value = -5267.500100000001

$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Range("H68").Value = 3569.6
$ws.Range("I68").Value = 1566.3334
$ws.Range("J68").Value = 4428.143
$ws.Range("K68").Value = 1566.3334
$ws.Range("L68").Value = 4428.143
$ws.Range("M68").Value = -817.3334
$ws.Range("N68").Value = -5926.143
# Row 71
$ws.Range("H71").Value = 3569.6
$ws.Range("I71").Value = 1566.3334
$ws.Range("J71").Value = 4428.143
$ws.Range("K71").Value = 7831.666999999999
$ws.Range("L71").Value = 22140.715
$ws.Range("M71").Value = -4087.666999999999
$ws.Range("N71").Value = -29628.715
# Row 93
$ws.Range("H93").Value = 6986
$ws.Range("I93").Value = 3488.6667
$ws.Range("K93").Value = 3488.6667
$ws.Range("M93").Value = -2240.6667

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 67249.25
$ws.Range("I62").Value = 102799
$ws.Range("J62").Value = 7999.6665
$ws.Range("K62").Value = 102799
$ws.Range("L62").Value = 7999.6665
$ws.Range("M62").Value = -102175
$ws.Range("N62").Value = -9247.666499999999
# Row 65
$ws.Range("H65").Value = 67249.25
$ws.Range("I65").Value = 102799
$ws.Range("J65").Value = 7999.6665
$ws.Range("K65").Value = 513995
$ws.Range("L65").Value = 39998.3325
$ws.Range("M65").Value = -510875
$ws.Range("N65").Value = -46238.3325
# Row 108
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()
# Row 113
$ws.Range("H113").Value = 735.4286
$ws.Range("I113").Value = 611.3889
$ws.Range("K113").Value = 1834.1667
$ws.Range("M113").Value = 335.8332999999998
